# Add test case for log out
# Target sheet: the first worksheet (internally named "Regestiration" in
# workbook.xml, but it is the sheet holding the Logout test cases -
# dimension A1:J3, TAWA_Logout_001 / TAWA_Logout_002 rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New row 3: "admin" logout test case -----------------------------
# A3 ("TAWA_Logout_002") already holds the right text - leave it alone.
$ws.Range("B3").Value = "Verify functionality of Logout link for admin"
$ws.Range("C3").Value = "Admin"
$ws.Range("C3").HorizontalAlignment = -4108   # xlCenter (matches C2's style)
$ws.Range("D3").Value = "Admin shall be logged in"
$ws.Range("E3").Value = '1-Click on "Logout" link in page Header in Admin page'

# --- Row 2 edits: tidy up / correct the existing "user" test case ----
$ws.Range("E2").Value = '1-Click on "Logout" link in page Header '

$ws.Range("F3").Value = "Admin shall be logged out and redirected to Home Page"
$ws.Range("F2").Value = "User shall be logged out and redirected to Home Page"

$ws.Range("G3").Value = "Crirical"
$ws.Range("G3").HorizontalAlignment = -4108   # xlCenter (matches G2's style)
$ws.Range("G2").Value = "Critical"

# --- View state: scroll right a bit and land the selection on F2 -----
[void]$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
[void]$ws.Range("F2").Select()
